$wb = $excel.ActiveWorkbook

# --- Sheet1: rename "Sheet1" -> "list-column" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "list-column"

# --- Add new sheet "two-row-header" right after the first sheet.
# Copying Sheet1 (instead of Worksheets.Add()) carries over the sheet's
# row-height/format defaults, then we wipe the copied cells and fill in
# the new table. ---
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Clear() | Out-Null
$ws2.Name = "two-row-header"

# Row 1: short column names
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "species"
$ws2.Range("C1").Value = "death"
$ws2.Range("D1").Value = "weight"

# Row 2: clarifying notes (written D,B,A,C to match the original
# shared-string insertion order of the authored workbook)
$ws2.Range("D2").Value = "(in grams)"
$ws2.Range("B2").Value = "(office supply type)"
$ws2.Range("A2").Value = "(at birth)"
$ws2.Range("C2").Value = "(date is approximate)"

# Row 3: actual data values
$ws2.Range("A3").Value = "Clippy"
$ws2.Range("B3").Value = "paperclip"
$ws2.Range("C3").Value = 39083
$ws2.Range("D3").Value = 0.9

# Give the date cell the same date format as Sheet1!B4, reusing the
# existing style record instead of fabricating a brand-new number format.
$ws1.Range("B4").Copy() | Out-Null
$ws2.Range("C3").PasteSpecial(-4122) | Out-Null

# Sheet1 keeps a plain selection on its data column, no longer the active tab
$ws1.Range("A2:A5").Select() | Out-Null

# Sheet2 becomes the active, visible tab with the header row selected
$ws2.Range("A1:D1").Select() | Out-Null
$ws2.Activate() | Out-Null
